$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-13 (12 rows), columns A-T (20 cols)
# Row order: ECs/ECs, ECs/FAPs, ECs/Resolving-Mac,
#            FAPs/ECs, FAPs/FAPs, FAPs/Resolving-Mac,
#            MuSCs/ECs, MuSCs/FAPs, MuSCs/Resolving-Mac,
#            Resolving-Mac/ECs, Resolving-Mac/FAPs, Resolving-Mac/Resolving-Mac
$data = @(
  @("ECs","Shank1","Sstr2","ECs",3,1,0.04115933333333333,0.123478,0.01731274819652447,0.01731274819652447,3,1,0.4513166666666666,1.35395,0.2951740315743233,0.2951740315743233,0.01857589312222222,0.1671830381,0.005110273682799223,0.005110273682799222),
  @("ECs","Shank1","Sstr2","FAPs",3,1,0.04115933333333333,0.123478,0.01731274819652447,0.01731274819652447,2,0.6666666666666666,1.049824666666667,3.149474,0.6866154126212269,0.6866154126212269,0.04321008339688889,0.388890750572,0.01188719974656405,0.01188719974656405),
  @("ECs","Shank1","Sstr2","Resolving-Mac",3,1,0.04115933333333333,0.123478,0.01731274819652447,0.01731274819652447,1,0.3333333333333333,0.02784366666666667,0.08353099999999999,0.01821055580444979,0.01821055580444979,0.001146026757555555,0.010314240818,0.0003152747671611964,0.0003152747671611963),
  @("FAPs","Shank1","Sstr2","ECs",3,1,1.977040333333333,5.931120999999999,0.8315975671465232,0.8315975671465232,3,1,0.4513166666666666,1.35395,0.2951740315743233,0.2951740315743233,0.8922712531055554,8.030441277949999,0.2454660065420383,0.2454660065420383),
  @("FAPs","Shank1","Sstr2","FAPs",3,1,1.977040333333333,5.931120999999999,0.8315975671465232,0.8315975671465232,2,0.6666666666666666,1.049824666666667,3.149474,0.6866154126212269,0.6866154126212269,2.075545708928222,18.679911380354,0.5709877067011185,0.5709877067011185),
  @("FAPs","Shank1","Sstr2","Resolving-Mac",3,1,1.977040333333333,5.931120999999999,0.8315975671465232,0.8315975671465232,1,0.3333333333333333,0.02784366666666667,0.08353099999999999,0.01821055580444979,0.01821055580444979,0.05504805202788888,0.4954324682509999,0.01514385390336645,0.01514385390336644),
  @("MuSCs","Shank1","Sstr2","ECs",3,1,0.322944,0.9688319999999999,0.135839132968911,0.135839132968911,3,1,0.4513166666666666,1.35395,0.2951740315743233,0.2951740315743233,0.1457500096,1.3117500864,0.04009618452399404,0.04009618452399404),
  @("MuSCs","Shank1","Sstr2","FAPs",3,1,0.322944,0.9688319999999999,0.135839132968911,0.135839132968911,2,0.6666666666666666,1.049824666666667,3.149474,0.6866154126212269,0.6866154126212269,0.3390345771519999,3.051311194368,0.09326924233355853,0.09326924233355854),
  @("MuSCs","Shank1","Sstr2","Resolving-Mac",3,1,0.322944,0.9688319999999999,0.135839132968911,0.135839132968911,1,0.3333333333333333,0.02784366666666667,0.08353099999999999,0.01821055580444979,0.01821055580444979,0.008991945087999999,0.08092750579199999,0.002473706111358429,0.002473706111358429),
  @("Resolving-Mac","Shank1","Sstr2","ECs",2,0.6666666666666666,0.03625666666666667,0.10877,0.01525055168804132,0.01525055168804132,3,1,0.4513166666666666,1.35395,0.2951740315743233,0.2951740315743233,0.01636323794444444,0.1472691415,0.004501566825491759,0.004501566825491759),
  @("Resolving-Mac","Shank1","Sstr2","FAPs",2,0.6666666666666666,0.03625666666666667,0.10877,0.01525055168804132,0.01525055168804132,2,0.6666666666666666,1.049824666666667,3.149474,0.6866154126212269,0.6866154126212269,0.03806314299777778,0.34256828698,0.01047126383998584,0.01047126383998584),
  @("Resolving-Mac","Shank1","Sstr2","Resolving-Mac",2,0.6666666666666666,0.03625666666666667,0.10877,0.01525055168804132,0.01525055168804132,1,0.3333333333333333,0.02784366666666667,0.08353099999999999,0.01821055580444979,0.01821055580444979,0.001009518541111111,0.00908566687,0.0002777210225637225,0.0002777210225637225)
)

# Clear the previously used range (rows 2-17) first so stale rows 14-17 are removed
$ws.Range("A2:T17").Clear()

# Write the new data rows starting at row 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $row[$j]
    }
}

$wb.Save()
